$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.538.28"

$ws.Range("D3").Value = "2.194.92"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "257.42"
$ws.Range("E5").Value = "  +1.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "84.30"
$ws.Range("E6").Value = "  +12.80%  "

$ws.Range("E7").Value = "  +0.91%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.592"
$ws.Range("E9").Value = "  +0.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.56"
$ws.Range("E10").Value = "  +8.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0915"
$ws.Range("E11").Value = "  +1.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.29"
$ws.Range("E12").Value = "  +6.17%  "

$ws.Range("E13").Value = "  +1.89%  "

$ws.Range("D14").Value = "2.526.38"
$ws.Range("E14").Value = "  -0.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.30"
$ws.Range("E15").Value = "  +0.08%  "

$ws.Range("D16").Value = "2.204.31"
$ws.Range("E16").Value = "  -0.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.779"
$ws.Range("E17").Value = "  -0.34%  "

$ws.Range("D18").Value = "43.470.24"
$ws.Range("E18").Value = "  +1.94%  "

$ws.Range("E19").Value = "  +0.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.52"
$ws.Range("E20").Value = "  -2.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.88"
$ws.Range("E21").Value = "  -0.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.32"
$ws.Range("E22").Value = "  +5.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.67"
$ws.Range("E23").Value = "  +0.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.91"
$ws.Range("E24").Value = "  -4.51%  "

$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.53"
$ws.Range("E26").Value = "  +4.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.58"
$ws.Range("E27").Value = "  +0.22%  "

$ws.Range("E28").Value = "  +2.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  +1.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.56"
$ws.Range("E30").Value = "  -1.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.16"
$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.31"
$ws.Range("E32").Value = "  +0.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0859"
$ws.Range("E33").Value = "  +1.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.27"
$ws.Range("E34").Value = "  +1.52%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0355"
$ws.Range("E37").Value = "  +1.91%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.44"
$ws.Range("E38").Value = "  +3.46%  "

$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.86"
$ws.Range("E39").Value = "  +3.59%  "

$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.25"
$ws.Range("E40").Value = "  -0.43%  "

$ws.Range("E41").Value = "  -0.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.02"
$ws.Range("E42").Value = "  +5.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.43"
$ws.Range("E43").Value = "  +3.68%  "

$ws.Range("E44").Value = "  +0.21%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0977"
$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "99.74"
$ws.Range("E46").Value = "  -1.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.28"
$ws.Range("E47").Value = "  -0.87%  "

$ws.Range("E48").Value = "  +3.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.10"
$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.432"
$ws.Range("E50").Value = "  -5.92%  "

$ws.Range("E51").Value = "  +3.10%  "
